$wb = $excel.ActiveWorkbook

$medium = $wb.Worksheets.Item("Medium")
$low    = $wb.Worksheets.Item("Low")
$high   = $wb.Worksheets.Item("High")

$medium.Range("G6").Value = 0.01627
$medium.Range("J6").Value = 0.017782
$medium.Range("H8").Value = 0.59
$medium.Range("K8").Value = 0.59

$low.Range("G6").Value = 0.01627
$low.Range("J6").Value = 0.017782
$low.Range("H8").Value = 0.59
$low.Range("K8").Value = 0.59

$high.Range("G6").Value = 0.01627
$high.Range("J6").Value = 0.017782
